$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates ---
$ws.Range("A1").Value = "Week 39 benefits"
$ws.Range("A2").Value = "Generated: 25/09/2025 13:07"
$ws.Rows(2).RowHeight = 6.75

# --- New header row (row 8) ---
$ws.Range("A8").Value = "Client Name"
$ws.Range("B8").Value = "SS Benefits"
$ws.Range("C8").Value = "Date Processed"

# --- New data row 9 (John Doe) ---
$ws.Range("A9").Value = "John Doe"
$ws.Range("B9").Value = 300
$ws.Range("C9").Value = "25/09/2025"

# --- New data row 10 (Jane Smith) ---
$ws.Range("A10").Value = "Jane Smith"
$ws.Range("B10").Value = 275.5
$ws.Range("C10").Value = "25/09/2025"

# --- Clear old client rows 12-14 (data removed) ---
$ws.Range("A12:D12").ClearContents()
$ws.Range("A13:D13").ClearContents()
$ws.Range("A14:D14").ClearContents()
